# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals.
# The "K" column (column G, header "K") values are recalculated and rewritten below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K (column G) value, per the recalculated s_vals.
$kValues = @{
    2  = 0
    3  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 4
    13 = 1
    14 = 1
    15 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
